$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-LabelCell($addr, $text) {
    $r = $ws.Range($addr)
    $r.Value = $text
    $r.Font.Bold = $true
    $r.WrapText = $false
    $r.VerticalAlignment = -4160
}
function Set-BlackCell($addr, $text) {
    $r = $ws.Range($addr)
    $r.Value = $text
    $r.Font.Bold = $false
    $r.WrapText = $true
    $r.VerticalAlignment = -4160
}
function Set-RedCell($addr, $text) {
    $r = $ws.Range($addr)
    $r.Value = $text
    $r.Font.Bold = $false
    $r.Font.Color = 255
    $r.WrapText = $true
    $r.VerticalAlignment = -4160
}

$ws.Range("A10:C26").Clear()

Set-LabelCell "A10" "Objetivos:"
Set-BlackCell "B10" "Apresentar os fundamentos da Termodinâmica e Física Estatística. Apresentar os diferentes formalismos da Física Estatística. Aplicação dos formalismos a modelos simples. Apresentar as aplicações."
Set-RedCell "C10" "Apresentar os fundamentos da Termodinâmica e Física Estatística. Apresentar os diferentes formalismos da Física Estatística. Aplicação dos formalismos a modelos simples. Apresentar as aplicações."
$ws.Rows.Item(10).RowHeight = 60

Set-LabelCell "A11" "Objectives:"
Set-BlackCell "B11" "This discipline aims to present the fundamentals of Thermodynamics and Statistical Physics. To present the different formalisms of Statistical Physics. Application of formalisms to some simple models. Applications."
Set-RedCell "C11" "This discipline aims to present the fundamentals of Thermodynamics and Statistical Physics. To present the different formalisms of Statistical Physics. Application of formalisms to some simple models. Applications."
$ws.Rows.Item(11).RowHeight = 60

Set-LabelCell "A12" "Docentes responsáveis:"
$ws.Rows.Item(12).AutoFit()

Set-BlackCell "B13" "1176388 - Luiz Tadeu Fernandes Eleno"
Set-RedCell "C13" "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Rows.Item(13).AutoFit()

Set-BlackCell "B14" "1643715 - Paulo Atsushi Suzuki"
Set-RedCell "C14" "1643715 - Paulo Atsushi Suzuki"
$ws.Rows.Item(14).AutoFit()

Set-LabelCell "A15" "Programa resumido:"
Set-BlackCell "B15" "Conceitos da termodinâmica. Formalismos da Física Estatística. Gás ideal clássico. Gases quânticos. Aplicações."
Set-RedCell "C15" "Conceitos da termodinâmica. Formalismos da Física Estatística. Gás ideal clássico. Gases quânticos. Aplicações."
$ws.Rows.Item(15).RowHeight = 60

Set-LabelCell "A16" "Short syllabus:"
Set-BlackCell "B16" "Concepts of thermodynamics. Formalisms of Statistical Physics. Classic ideal gas. Quantum gases. Applications."
Set-RedCell "C16" "Concepts of thermodynamics. Formalisms of Statistical Physics. Classic ideal gas. Quantum gases. Applications."
$ws.Rows.Item(16).RowHeight = 60

Set-LabelCell "A17" "Programa:"
Set-BlackCell "B17" "Sistemas macroscópicos e microscópicos. Postulados da termodinâmica. Equação fundamental. Equações de estado. Equação de Euler. Relação de Gibbs-Duhem. Equilíbrio termodinâmico. Derivadas Termodinâmicas. Potenciais Termodinâmicos: Helmholtz, Gibbs, Grande Canônico e Entalpia. Relações de Maxwell. Diagrama de Born. Redução de derivadas termodinâmicas. Formalismo microcanônico. Equação de Boltzmann. Exemplo: Modelo de Einstein de sólido cristalino. Formalismo canônico. Exemplo: Gás ideal clássico. Distribuição de Maxwell-Boltzmann. Formalismo grande canônico. Gases quânticos: férmions e bósons. Estatística de Bose-Einstein. Estatística de Fermi-Dirac. Exemplos: gás de elétrons e gás de fótons. Estatística de Planck.Aplicações: supercondutividade. Gás de elétrons em semicondutores. superfluidez no hélio líquido."
Set-RedCell "C17" "Sistemas macroscópicos e microscópicos. Postulados da termodinâmica. Equação fundamental. Equações de estado. Equação de Euler. Relação de Gibbs-Duhem. Equilíbrio termodinâmico. Derivadas Termodinâmicas. Potenciais Termodinâmicos: Helmholtz, Gibbs, Grande Canônico e Entalpia. Relações de Maxwell. Diagrama de Born. Redução de derivadas termodinâmicas. Formalismo microcanônico. Equação de Boltzmann. Exemplo: Modelo de Einstein de sólido cristalino. Formalismo canônico. Exemplo: Gás ideal clássico. Distribuição de Maxwell-Boltzmann. Formalismo grande canônico. Gases quânticos: férmions e bósons. Estatística de Bose-Einstein. Estatística de Fermi-Dirac. Exemplos: gás de elétrons e gás de fótons. Estatística de Planck.Aplicações: supercondutividade. Gás de elétrons em semicondutores. superfluidez no hélio líquido."
$ws.Rows.Item(17).RowHeight = 120

Set-LabelCell "A18" "Syllabus:"
Set-BlackCell "B18" "Macroscopic and microscopic systems. Postulates of thermodynamics. Fundamental equation. Equations of state. Thermodynamic equilibrium. Thermodynamic derivatives. Thermodynamic potentials. Maxwell relations. Born diagram. Reduction of thermodynamic derivatives.Microcanonical formalism. Boltzmann equation. Einstein model of a crystalline solid. Canonical formalism. Example: ideal classical gas. Maxwell-Boltzmann distribution. Grand canonical formalism. Quantum gases. Fermions and bosons. Bose-Einstein distribution. Fermi-Dirac distribution. Examples: electron gas and photon gas. Planck distribution.Applications: superconductivity, electron gas in semiconductor, superfluidity of the liquid helium."
Set-RedCell "C18" "Macroscopic and microscopic systems. Postulates of thermodynamics. Fundamental equation. Equations of state. Thermodynamic equilibrium. Thermodynamic derivatives. Thermodynamic potentials. Maxwell relations. Born diagram. Reduction of thermodynamic derivatives.Microcanonical formalism. Boltzmann equation. Einstein model of a crystalline solid. Canonical formalism. Example: ideal classical gas. Maxwell-Boltzmann distribution. Grand canonical formalism. Quantum gases. Fermions and bosons. Bose-Einstein distribution. Fermi-Dirac distribution. Examples: electron gas and photon gas. Planck distribution.Applications: superconductivity, electron gas in semiconductor, superfluidity of the liquid helium."
$ws.Rows.Item(18).RowHeight = 120

Set-LabelCell "A19" "Avaliação:"
$ws.Rows.Item(19).AutoFit()

Set-LabelCell "A20" "Método:"
Set-BlackCell "B20" "Aulas expositivas, seminários e exercícios comentados."
Set-RedCell "C20" "Aulas expositivas, seminários e exercícios comentados."
$ws.Rows.Item(20).RowHeight = 60

Set-LabelCell "A21" "Critério:"
Set-BlackCell "B21" "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."
Set-RedCell "C21" "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."
$ws.Rows.Item(21).RowHeight = 60

Set-LabelCell "A22" "Norma de recuperação:"
Set-BlackCell "B22" "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
Set-RedCell "C22" "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Rows.Item(22).RowHeight = 60

Set-LabelCell "A23" "Bibliografia:"
Set-BlackCell "B23" "CALLEN, H.B., Thermodynamics and an introduction to thermostatistics, John Wiley & Sons, New York, 1985. SALINAS, S. R. A., Introdução à Física Estatística, Edusp, São Paulo, 1999. CASQUILHO J.P. e TEIXEIRA P.I.C, Introdução à Física Estatística, Editora Livraria da Física, São Paulo, 2012. DALARSSON, N. DALARSSON, M. GOLUBOVIC, L.  Introductory Statistical Thermodynamics. Academic Press, 2011."
Set-RedCell "C23" "CALLEN, H.B., Thermodynamics and an introduction to thermostatistics, John Wiley & Sons, New York, 1985. SALINAS, S. R. A., Introdução à Física Estatística, Edusp, São Paulo, 1999. CASQUILHO J.P. e TEIXEIRA P.I.C, Introdução à Física Estatística, Editora Livraria da Física, São Paulo, 2012. DALARSSON, N. DALARSSON, M. GOLUBOVIC, L.  Introductory Statistical Thermodynamics. Academic Press, 2011."
$ws.Rows.Item(23).RowHeight = 120

Set-LabelCell "A24" "Requisitos:"
$ws.Rows.Item(24).AutoFit()

Set-BlackCell "B25" "LOB1019 -  Física II  (Requisito)`n"
Set-RedCell "C25" "LOB1019 -  Física II  (Requisito)`n"
$ws.Rows.Item(25).RowHeight = 30

Set-BlackCell "B26" "LOB1052 -  Cálculo III  (Requisito)`n"
Set-RedCell "C26" "LOB1052 -  Cálculo III  (Requisito)`n"
$ws.Rows.Item(26).RowHeight = 30
